$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.787.45"
$ws.Range("E2").Value = "  +4.60%  "
$ws.Range("D3").Value = "2.259.03"
$ws.Range("E3").Value = "  +4.11%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'249.44"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D7").Value = "'71.56"
$ws.Range("E7").Value = "  +8.24%  "
$ws.Range("D9").Value = "'0.658"
$ws.Range("E9").Value = "  +16.86%  "
$ws.Range("D10").Value = "'38.87"
$ws.Range("E10").Value = "  +9.52%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0970"
$ws.Range("E11").Value = "  +4.72%  "
$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").Value = "'59.57"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").Value = "'7.45"
$ws.Range("E13").Value = "  +8.62%  "
$ws.Range("D14").Value = "'0.105"
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").Value = "2.590.06"
$ws.Range("E15").Value = "  +3.89%  "
$ws.Range("D16").Value = "'14.88"
$ws.Range("E16").Value = "  +4.63%  "
$ws.Range("D17").Value = "'0.883"
$ws.Range("E17").Value = "  +3.62%  "
$ws.Range("D18").Value = "2.264.86"
$ws.Range("E18").Value = "  +4.28%  "
$ws.Range("D19").Value = "42.730.10"
$ws.Range("E19").Value = "  +4.48%  "
$ws.Range("E20").Value = "  +6.12%  "
$ws.Range("D21").Value = "'6.31"
$ws.Range("E21").Value = "  +3.92%  "
$ws.Range("D22").Value = "'73.09"
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("D23").Value = "'235.78"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").Value = "'2.09"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  +7.21%  "
$ws.Range("D26").Value = "'11.47"
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("D29").Value = "'3.67"
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("E30").Value = "  +5.90%  "
$ws.Range("D31").Value = "'167.69"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").Value = "'20.95"
$ws.Range("E32").Value = "  +4.04%  "
$ws.Range("D33").Value = "'6.52"
$ws.Range("E33").Value = "  +15.29%  "
$ws.Range("D34").Value = "'0.127"
$ws.Range("E34").Value = "  +5.27%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'31.49"
$ws.Range("E35").Value = "  +28.37%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.0798"
$ws.Range("E36").Value = "  +6.79%  "
$ws.Range("E37").Value = "  +4.08%  "
$ws.Range("D38").Value = "'4.45"
$ws.Range("E38").Value = "  +12.80%  "
$ws.Range("D39").Value = "'4.75"
$ws.Range("E39").Value = "  +4.38%  "
$ws.Range("E40").Value = "  +4.70%  "
$ws.Range("E41").Value = "  +6.32%  "
$ws.Range("D42").Value = "'12.62"
$ws.Range("E42").Value = "  +10.91%  "
$ws.Range("E43").Value = "  +6.61%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'9.14"
$ws.Range("E44").Value = "  +8.85%  "
$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").Value = "'62.17"
$ws.Range("E45").Value = "  +3.41%  "
$ws.Range("E46").Value = "  +5.19%  "
$ws.Range("D47").Value = "'4.87"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("E48").Value = "  +3.60%  "
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("E50").Value = "  +3.33%  "
$ws.Range("E51").Value = "  +4.34%  "
